$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated 2D training schedules values (columns C-H for rows 2-6)
$values = @{
    2 = @{ C=6;  D=5; E=3; F=3; G=-3; H=34 }
    3 = @{ C=8;  D=2; E=3; F=1; G=-5; H=56 }
    4 = @{ C=9;  D=7; E=7; F=4; G=-2; H=23 }
    5 = @{ C=5;  D=5; E=1; F=2; G=-4; H=45 }
    6 = @{ C=7;  D=7; E=6; F=5; G=-1; H=12 }
}

# Column B values also changed for some rows
$bValues = @{
    3 = 1
    4 = 3
    6 = 2
}

foreach ($row in $bValues.Keys) {
    $ws.Range("B$row").Value = $bValues[$row]
}

foreach ($row in $values.Keys) {
    $rowData = $values[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
